# ------------------------------------------------------------------
# Applies the "Updating system state machine, cleanup, EEPROM store
# update" edit to List1 (sheet1): new input values, a highlighted
# (orange/green) boxed summary table with borders around A4:D12, a
# couple of new rows, renamed TO_att/TO_f0 labels (now Greek symbols
# TO_ξ / TO_ω0), a new BEMF_DQ title and a new "TO" label at F8, plus
# updated formulas in A12/B12 and the resulting downstream values.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# ---- colors (Excel .Interior.Color / .Font.Color use BGR integers) ----
$orange = 49407   # RGB(255,192,0)  -> FFC000
$green  = 5287936 # RGB(0,176,80)   -> 00B050

# ---- header row (row 1) : font only changes (reuses existing header font) ----
$ws.Range("H1:L1").Font.Name = "Calibri"

# ---- row 2: new input values ----
$ws.Range("A2").Value = 0.0023
$ws.Range("B2").Value = 0.0023
$ws.Range("C2").Value = 0.02
$ws.Range("G2").Value = 20000

# row 2 fills: A2:F2 and J2 get the orange fill (input values w/ sci fmt)
$ws.Range("A2:F2").Interior.Color = $orange
$ws.Range("J2").Interior.Color = $orange
# G2, K2, L2 orange too, but plain General format (no sci notation)
$ws.Range("G2").Interior.Color = $orange
$ws.Range("K2").Interior.Color = $orange
$ws.Range("L2").Interior.Color = $orange
# H2, I2 (derived omega values) get the green fill
$ws.Range("H2:I2").Interior.Color = $green

# ---- row 3: thin spacer row below the input row ----
$ws.Rows.Item(3).RowHeight = 15.75

# ---- row 5: damping ratio / natural frequency labels (now Greek letters) ----
# (B5 -> TO_ω0 is created before A5 -> TO_ξ so the shared-string table
# ends up in the same append order as the authored workbook.)
$ws.Range("B5").Value = "TO_" + [char]0x3C9 + "0"    # TO_ω0
$chB5w = $ws.Range("B5").Characters(4, 1)
$chB5w.Font.Name = "Calibri"
$chB5sub = $ws.Range("B5").Characters(5, 1)
$chB5sub.Font.Subscript = $true

$ws.Range("A5").Value = "TO_" + [char]0x3BE          # TO_ξ
$chA5 = $ws.Range("A5").Characters(4, 1)
$chA5.Font.Name = "Calibri"

$ws.Range("C5").Value = $null
$ws.Range("D5").Value = $null
$ws.Rows.Item(5).RowHeight = 18

# ---- Box header: A4 becomes the "BEMF_DQ" title bar (orange-less, bordered) ----
$ws.Range("A4").Value = "BEMF_DQ"
$ws.Range("B4").Value = $null
$ws.Range("C4").Value = $null
$ws.Range("D4").Value = $null

# ---- row 6: A6 0.8 -> 1, B6 unchanged; pad C6/D6 inside the box ----
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 40
$ws.Range("C6").Value = $null
$ws.Range("D6").Value = $null
$ws.Range("A6:B6").Interior.Color = $orange

# ---- row 7: blank spacer row inside the box ----
$ws.Range("A7").Value = $null
$ws.Range("B7").Value = $null
$ws.Range("C7").Value = $null
$ws.Range("D7").Value = $null

# ---- row 8: Ufrac/Efrac/Wifrac/Ifrac header (unchanged text), plus new F8 "TO" ----
$ws.Range("F8").Value = "TO"

# ---- row 9: results row - green fill, new number format ----
$ws.Range("A9").NumberFormat = "0.00E+00"
$ws.Range("B9:D9").NumberFormat = "0.000000"
$ws.Range("A9:D9").Interior.Color = $green

# ---- row 10: blank spacer row inside the box ----
$ws.Range("A10").Value = $null
$ws.Range("B10").Value = $null
$ws.Range("C10").Value = $null
$ws.Range("D10").Value = $null

# ---- row 11: TO_Kp/TO_Ki/TO_Th header (unchanged text); pad D11 ----
$ws.Range("D11").Value = $null

# ---- row 12: formulas updated to factor in A2/C2; thick bottom border row ----
$ws.Range("A12").Formula = "=(2*A6*2*PI()*B6*A2)-C2"
$ws.Range("B12").Formula = "=((2*PI()*B6)^2)*A2"
$ws.Range("C12").Value = $null
$ws.Range("D12").Value = $null
$ws.Range("A12:B12").NumberFormat = "0.00E+00"
$ws.Range("A12:D12").Interior.Color = $green
$ws.Rows.Item(12).RowHeight = 15.75

# ---- row 15: derived values now use the green fill / plain number format ----
$ws.Range("A15:C15").Interior.Color = $green

# ------------------------------------------------------------------
# Borders: a medium box drawn around A4:D12
# ------------------------------------------------------------------
$medium = -4138

# outer edges
$ws.Range("A4:D4").Borders.Item(8).Weight = $medium   # xlEdgeTop
$ws.Range("A12:D12").Borders.Item(9).Weight = $medium # xlEdgeBottom
$ws.Range("A4:A12").Borders.Item(7).Weight = $medium  # xlEdgeLeft
$ws.Range("D4:D12").Borders.Item(10).Weight = $medium # xlEdgeRight

# ------------------------------------------------------------------
# Selection
# ------------------------------------------------------------------
$ws.Range("D9").Select()

Write-Host "Edit applied"
